$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Helper to set a cell value, preserving "numeric-looking" strings (like "00")
# as literal text instead of letting them be auto-converted to numbers.
function Set-TextValue($cell, $value) {
    $isNumeric = $value -match '^-?\d+(\.\d+)?$'

    if ($isNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

# Full target data for rows 4-13, columns B-F
$rows = @(
    @{ Row = 4;  B = "SingleUseId1";  C = "Large";   D = "Left";   E = "LTR"; F = "<hour>" },
    @{ Row = 5;  B = "SingleUseId2";  C = "Large";   D = "Left";   E = "LTR"; F = "00" },
    @{ Row = 6;  B = "SingleUseId3";  C = "Large";   D = "Left";   E = "LTR"; F = "<minute>" },
    @{ Row = 7;  B = "SingleUseId4";  C = "Large";   D = "Left";   E = "LTR"; F = "00" },
    @{ Row = 8;  B = "SingleUseId5";  C = "Default"; D = "Center"; E = "LTR"; F = "Set" },
    @{ Row = 9;  B = "SingleUseId6";  C = "Default"; D = "Center"; E = "LTR"; F = "Set" },
    @{ Row = 10; B = "SingleUseId7";  C = "Default"; D = "Center"; E = "LTR"; F = "Clock" },
    @{ Row = 11; B = "SingleUseId8";  C = "Large";   D = "Left";   E = "LTR"; F = "<hour>:<minute>" },
    @{ Row = 12; B = "SingleUseId9";  C = "Large";   D = "Left";   E = "LTR"; F = "00" },
    @{ Row = 13; B = "SingleUseId10"; C = "Large";   D = "Left";   E = "LTR"; F = "00" }
)

foreach ($r in $rows) {
    Set-TextValue $ws.Cells.Item($r.Row, 2) $r.B
    Set-TextValue $ws.Cells.Item($r.Row, 3) $r.C
    Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D
    Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E
    Set-TextValue $ws.Cells.Item($r.Row, 6) $r.F
}
